$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.404.92"
$ws.Range("E2").Value = "  +0.32%  "

# Row 3
$ws.Range("D3").Value = "1.878.17"
$ws.Range("E3").Value = "  +0.16%  "

# Row 4
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7164"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.87%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "243.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.83%  "

# Row 7
$ws.Range("E7").Value = "  +0.08%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07981"
$ws.Range("D8").Style = "Normal"

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3147"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.71%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.93"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.14%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08082"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.90%  "

# Row 12
$ws.Range("D12").Value = "1.879.43"
$ws.Range("E12").Value = "  +0.18%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "94.72"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.84%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.229"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.26%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7079"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.27%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.403"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.25%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008450"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.14%  "

# Row 18
$ws.Range("D18").Value = "29.409.57"
$ws.Range("E18").Value = "  +0.34%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "253.26"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.25%  "

# Row 20
$ws.Range("E20").Value = "  +0.90%  "

# Row 21
$ws.Range("D21").Value = "2.135.27"
$ws.Range("E21").Value = "  +0.41%  "

# Row 22
$ws.Range("E22").Value = "  +0.11%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.680"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.46%  "

# Row 24
$ws.Range("E24").Value = "  +0.03%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1578"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.84%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.071"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.20%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.40%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.97"
$ws.Range("D28").Style = "Normal"

# Row 29
$ws.Range("E29").Value = "  +0.28%  "

# Row 30
$ws.Range("E30").Value = "  -0.04%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.318"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.65%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.223"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.47%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05309"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.87%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.944"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.01%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7583"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.07%  "

# Row 36
$ws.Range("E36").Value = "  -0.11%  "

# Row 37
$ws.Range("E37").Value = "  +0.28%  "

# Row 38
$ws.Range("E38").Value = "  +0.15%  "

# Row 39
$ws.Range("D39").Value = "1.275.93"
$ws.Range("E39").Value = "  -1.38%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.758"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.74%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.414"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.57%  "

# Row 42
$ws.Range("E42").Value = "  +1.38%  "

# Row 43
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "74.19"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.21%  "

# Row 44
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "111.64"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.59%  "

# Row 45
$ws.Range("E45").Value = "  +0.07%  "

# Row 46
$ws.Range("E46").Value = "  -0.56%  "

# Row 47
$ws.Range("D47").Value = "2.029.95"
$ws.Range("E47").Value = "  +0.42%  "

# Row 48
$ws.Range("E48").Value = "  +0.25%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.5207"
$ws.Range("D49").Style = "Normal"

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.524"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.78%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4344"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.32%  "
